# Generate Report for handoff
# Renames the handed-off markdown source file, records a new handoff
# attempt for it, and adds a "Handoff failed" row for a second file
# (which pushes the ".localization-config" bookkeeping row down by one).

$wb = $excel.ActiveWorkbook

$newFile = "e2212604-b45a-4ebc-af33-2c6f14975545"
$failedFile = "aca59f08-a9bf-453d-83b4-6b3f1e5f661b"
$newHash = "c724ce3357910043431093a707bfc7259c968a11"

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/dfbbd77ff8a20f1b2d0ef20bfb987fef86605a11"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dcba4d004d8c7c9ad9fe46079121668c1b92838b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb6fa8012ad40ead91e225f652ad3ebed10dab30/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

$mdDisplayNew = $newFile + ".md"
$mdDisplayFailed = $failedFile + ".md"
$configDisplay = ".localization-config"

$zhXlfDisplay = $newFile + "." + $newHash + ".zh-cn.xlf"
$deXlfDisplay = $newFile + "." + $newHash + ".de-de.xlf"

$mdUrlNew = $repoBase + "/e2e/" + $mdDisplayNew
$mdUrlFailed = $repoBase + "/e2e/" + $mdDisplayFailed
$configUrl = $repoBase + "/" + $configDisplay
$zhXlfUrl = $zhHandoffBase + "/" + $zhXlfDisplay
$deXlfUrl = $deHandoffBase + "/" + $deXlfDisplay

$newHandoffTimeZh = "2016-01-08 10:35:47"
$newHandoffTimeDe = "2016-01-08 10:35:57"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B3").Value = "Handoff failed"
$ws1.Range("C3").Value = "Handoff failed"
$ws1.Range("B4").Value = "Not localized"
$ws1.Range("C4").Value = "Not localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdUrlNew, [Type]::Missing, [Type]::Missing, $mdDisplayNew) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdUrlFailed, [Type]::Missing, [Type]::Missing, $mdDisplayFailed) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, [Type]::Missing, [Type]::Missing, $configDisplay) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("D2").Value = $newHandoffTimeZh
$ws2.Range("B3").Value = "Handoff failed"

$ws2.Range("A4").Value = $configDisplay
$ws2.Range("B4").Value = "Not localized"
$ws2.Range("D4").Value = $epoch
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrlNew, [Type]::Missing, [Type]::Missing, $mdDisplayNew) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfDisplay) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdUrlFailed, [Type]::Missing, [Type]::Missing, $mdDisplayFailed) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, [Type]::Missing, [Type]::Missing, $configDisplay) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("D2").Value = $newHandoffTimeDe
$ws3.Range("B3").Value = "Handoff failed"

$ws3.Range("A4").Value = $configDisplay
$ws3.Range("B4").Value = "Not localized"
$ws3.Range("D4").Value = $epoch
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrlNew, [Type]::Missing, [Type]::Missing, $mdDisplayNew) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfDisplay) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdUrlFailed, [Type]::Missing, [Type]::Missing, $mdDisplayFailed) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, [Type]::Missing, [Type]::Missing, $configDisplay) | Out-Null
